# Auto-generated Excel COM-interop edit script
# Applies the cell-level numeric updates described by the target diff
# against Sheets/Phantom_Profits.xlsx (all 8 job sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook


# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 5185.4443
$ws.Range("I32").Value = 5099.5
$ws.Range("J32").Value = 5254.2
$ws.Range("K32").Value = 5099.5
$ws.Range("L32").Value = 5254.2
$ws.Range("M32").Value = -4773.5
$ws.Range("N32").Value = -5906.2
$ws.Range("H40").Value = 3649.7
$ws.Range("I40").Value = 3299.4
$ws.Range("K40").Value = 3299.4
$ws.Range("M40").Value = -3124.4
$ws.Range("H74").Value = 5993
$ws.Range("I74").Value = 5993
$ws.Range("K74").Value = 5993
$ws.Range("M74").Value = -5057
$ws.Range("H77").Value = 5993
$ws.Range("I77").Value = 5993
$ws.Range("K77").Value = 29965
$ws.Range("M77").Value = -25285
$ws.Range("H98").Value = 365
$ws.Range("I98").Value = 365.8
$ws.Range("J98").Value = 361
$ws.Range("K98").Value = 365.8
$ws.Range("L98").Value = 361
$ws.Range("M98").Value = 1132.2
$ws.Range("N98").Value = -3357
$ws.Range("H107").Value = 1334
$ws.Range("I107").Value = 1271.1111
$ws.Range("J107").Value = 1900
$ws.Range("K107").Value = 1271.1111
$ws.Range("L107").Value = 1900
$ws.Range("M107").Value = 648.8888999999999
$ws.Range("N107").Value = -5740
$ws.Range("H111").Value = 1894.5
$ws.Range("I111").Value = 1894.5
$ws.Range("K111").Value = 5683.5
$ws.Range("M111").Value = -2616.5
$ws.Range("H113").Value = 5265.636
$ws.Range("J113").Value = 7599.25
$ws.Range("L113").Value = 7599.25
$ws.Range("N113").Value = -14107.25
$ws.Range("H122").Value = 365
$ws.Range("I122").Value = 365.8
$ws.Range("J122").Value = 361
$ws.Range("K122").Value = 1097.4
$ws.Range("L122").Value = 1083
$ws.Range("M122").Value = 1352.6
$ws.Range("N122").Value = -5983
$ws.Range("H135").Value = 2041
$ws.Range("I135").Value = 562.25
$ws.Range("J135").Value = 4998.5
$ws.Range("K135").Value = 5060.25
$ws.Range("L135").Value = 44986.5
$ws.Range("M135").Value = -2525.25
$ws.Range("N135").Value = -50056.5
$ws.Range("H137").Value = 1794.4
$ws.Range("I137").Value = 1725.5385
$ws.Range("K137").Value = 5176.6155
$ws.Range("M137").Value = -2626.6155

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2092
$ws.Range("I88").Value = 1373.5
$ws.Range("J88").Value = 2666.8
$ws.Range("K88").Value = 1373.5
$ws.Range("L88").Value = 2666.8
$ws.Range("M88").Value = -967.5
$ws.Range("N88").Value = -3478.8
$ws.Range("H91").Value = 2092
$ws.Range("I91").Value = 1373.5
$ws.Range("J91").Value = 2666.8
$ws.Range("K91").Value = 1373.5
$ws.Range("L91").Value = 2666.8
$ws.Range("M91").Value = 30.5
$ws.Range("N91").Value = -5474.8
$ws.Range("H122").Value = 1901.5
$ws.Range("I122").Value = 1901.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5704.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3254.5
$ws.Range("N122").ClearContents()
$ws.Range("H125").Value = 100500
$ws.Range("J125").Value = 100500
$ws.Range("L125").Value = 100500
$ws.Range("N125").Value = -110340
$ws.Range("H132").Value = 2531.25
$ws.Range("I132").Value = 1480.1666
$ws.Range("J132").Value = 3582.3333
$ws.Range("K132").Value = 4440.4998
$ws.Range("L132").Value = 10746.9999
$ws.Range("M132").Value = -1910.4998
$ws.Range("N132").Value = -15806.9999
$ws.Range("H135").Value = 71666.664
$ws.Range("J135").Value = 71666.664
$ws.Range("L135").Value = 71666.664
$ws.Range("N135").Value = -81806.664

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H20").Value = 1162.25
$ws.Range("I20").Value = 1162.25
$ws.Range("K20").Value = 1162.25
$ws.Range("M20").Value = -915.25

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5083.3335
$ws.Range("I31").Value = 4166.6665
$ws.Range("J31").Value = 6000
$ws.Range("K31").Value = 4166.6665
$ws.Range("L31").Value = 6000
$ws.Range("M31").Value = -3871.6665
$ws.Range("N31").Value = -6590
$ws.Range("H34").Value = 5083.3335
$ws.Range("I34").Value = 4166.6665
$ws.Range("J34").Value = 6000
$ws.Range("K34").Value = 4166.6665
$ws.Range("L34").Value = 6000
$ws.Range("M34").Value = -3964.6665
$ws.Range("N34").Value = -6404
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 183.5
$ws.Range("I12").Value = 202.66667
$ws.Range("K12").Value = 608.00001
$ws.Range("M12").Value = -435.00001
$ws.Range("H23").Value = 321.5
$ws.Range("I23").Value = 247.25
$ws.Range("J23").Value = 470
$ws.Range("K23").Value = 741.75
$ws.Range("L23").Value = 1410
$ws.Range("M23").Value = -506.75
$ws.Range("N23").Value = -1880
$ws.Range("H98").Value = 392.2
$ws.Range("J98").Value = 416.25
$ws.Range("L98").Value = 1248.75
$ws.Range("N98").Value = -4244.75
$ws.Range("H124").Value = 8578.4
$ws.Range("I124").Value = 4949.5
$ws.Range("J124").Value = 10997.667
$ws.Range("K124").Value = 14848.5
$ws.Range("L124").Value = 32993.001
$ws.Range("M124").Value = -9938.5
$ws.Range("N124").Value = -42813.001
$ws.Range("H131").Value = 1683.25
$ws.Range("I131").Value = 1429.4445
$ws.Range("J131").Value = 1890.909
$ws.Range("K131").Value = 4288.333500000001
$ws.Range("L131").Value = 5672.727000000001
$ws.Range("M131").Value = 751.6664999999994
$ws.Range("N131").Value = -15752.727

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4907.8335
$ws.Range("I102").Value = 4907.8335
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 4907.8335
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -3285.8335
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 2881.2856
$ws.Range("I122").Value = 2309.7778
$ws.Range("J122").Value = 3910
$ws.Range("K122").Value = 6929.3334
$ws.Range("L122").Value = 11730
$ws.Range("M122").Value = -4479.3334
$ws.Range("N122").Value = -16630
$ws.Range("H126").Value = 2772.6667
$ws.Range("I126").Value = 2772.6667
$ws.Range("K126").Value = 8318.000100000001
$ws.Range("M126").Value = -5848.000100000001
$ws.Range("H132").Value = 3545.0476
$ws.Range("I132").Value = 3550.4707
$ws.Range("K132").Value = 10651.4121
$ws.Range("M132").Value = -8121.4121

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1901.7693
$ws.Range("I7").Value = 1652.3
$ws.Range("J7").Value = 2733.3333
$ws.Range("K7").Value = 1652.3
$ws.Range("L7").Value = 2733.3333
$ws.Range("M7").Value = -1540.3
$ws.Range("N7").Value = -2957.3333
$ws.Range("H40").Value = 7937.385
$ws.Range("I40").Value = 4432.1665
$ws.Range("K40").Value = 4432.1665
$ws.Range("M40").Value = -4296.1665
$ws.Range("H46").Value = 2576.7778
$ws.Range("I46").Value = 2115.1667
$ws.Range("K46").Value = 2115.1667
$ws.Range("M46").Value = -1927.1667
$ws.Range("H93").Value = 968.6
$ws.Range("I93").Value = 968.6
$ws.Range("K93").Value = 968.6
$ws.Range("M93").Value = 279.4
$ws.Range("H100").Value = 2310.25
$ws.Range("I100").Value = 2310.25
$ws.Range("K100").Value = 2310.25
$ws.Range("M100").Value = -1769.25
$ws.Range("H122").Value = 14624.5
$ws.Range("I122").Value = 14624.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 43873.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -41423.5
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 1901.7693
$ws.Range("I126").Value = 1652.3
$ws.Range("J126").Value = 2733.3333
$ws.Range("K126").Value = 4956.9
$ws.Range("L126").Value = 8199.999899999999
$ws.Range("M126").Value = -2486.9
$ws.Range("N126").Value = -13139.9999
$ws.Range("H132").Value = 1989.7391
$ws.Range("I132").Value = 1940.421
$ws.Range("J132").Value = 2224
$ws.Range("K132").Value = 5821.263
$ws.Range("L132").Value = 6672
$ws.Range("M132").Value = -3291.263
$ws.Range("N132").Value = -11732

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4332.75
$ws.Range("I81").Value = 3932.8333
$ws.Range("J81").Value = 4732.6665
$ws.Range("K81").Value = 7865.6666
$ws.Range("L81").Value = 9465.333000000001
$ws.Range("M81").Value = -6804.6666
$ws.Range("N81").Value = -11587.333
$ws.Range("H84").Value = 4332.75
$ws.Range("I84").Value = 3932.8333
$ws.Range("J84").Value = 4732.6665
$ws.Range("K84").Value = 39328.333
$ws.Range("L84").Value = 47326.665
$ws.Range("M84").Value = -34024.333
$ws.Range("N84").Value = -57934.665
$ws.Range("H96").Value = 2599.8333
$ws.Range("I96").Value = 2000
$ws.Range("J96").Value = 2899.75
$ws.Range("K96").Value = 2000
$ws.Range("L96").Value = 2899.75
$ws.Range("M96").Value = -627
$ws.Range("N96").Value = -5645.75
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H107").Value = 1826
$ws.Range("I107").Value = 805.5
$ws.Range("K107").Value = 2416.5
$ws.Range("M107").Value = -496.5
$ws.Range("H122").Value = 5662
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 3088.2856
$ws.Range("I126").Value = 3088.2856
$ws.Range("K126").Value = 9264.856800000001
$ws.Range("M126").Value = -6794.856800000001
$ws.Range("H136").Value = 3308.125
$ws.Range("I136").Value = 3154.25
$ws.Range("K136").Value = 9462.75
$ws.Range("M136").Value = -6912.75
